$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Data de Conclusão" (completion date) values in column B (rows 2-9)
# with their new serial-date values, keeping the existing date number format.
$ws.Range("B2").Value = 43784
$ws.Range("B3").Value = 43785
$ws.Range("B4").Value = 43789
$ws.Range("B5").Value = 43798
$ws.Range("B6").Value = 43806
$ws.Range("B7").Value = 43820
$ws.Range("B8").Value = 43834
$ws.Range("B9").Value = 43839

# New row 10 appears with B10 carrying the same style as the rows above it,
# but with no value (empty cell) - copy the style from B9 without the value.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = $null

# Move the active selection to D13 (matches the post-edit sheetView selection).
$ws.Range("D13").Select() | Out-Null
